$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = "pants for men sport"
    2 = "mens pouch leggings"
    3 = "spandex leggings boys"
    4 = "rodilleras de basketball"
    5 = "knee sleeves wrestling"
    6 = "knee pads workout"
    7 = "baseball catcher leg guards adult"
    8 = "5 inch seam shorts men"
    9 = "youth softball pants for girls"
    10 = "mens compression workout pants"
    11 = "compression shorts for men"
    12 = "work knee pads for men"
    13 = "padded shorts snowboarding"
    14 = "youth girls yoga pants"
    15 = "youth hockey padded shorts"
    16 = "knee protector sports"
    17 = "kneeling pad construction"
    18 = "knee pad for work"
    19 = "basketballs under"
    20 = "boys baseball pants size 6"
    21 = "sliding short"
    22 = "black knee pads volleyball girls"
    23 = "mens 3/4 pants"
    24 = "football leg sleeves for men"
    25 = "black basketball shorts men"
    26 = "compression football shorts"
    27 = "running leggings men"
    28 = "football leggings boys"
    29 = "baseball leg guards"
    30 = "protective basketball"
    31 = "mountain bike knee pads"
    32 = "weightlifting shorts men"
    33 = "shorts for men basketball"
    34 = "wrestling shorts for boys"
    35 = "baseball items for men"
    36 = "knee pads for biking"
    37 = "eva foam knee pads"
    38 = "mens compression running tights"
    39 = "womens softball pants black"
    40 = "waist guard"
    41 = "bump pads"
    42 = "mens fitness pants"
    43 = "cycling pants for men padded"
    44 = "knee pads work"
    45 = "youth football girdle"
    46 = "bjj knee sleeves"
    47 = "volleyball kneepads black"
    48 = "non slip knee pads"
    49 = "short baseball"
    50 = "knee sleeves for wrestling"
    51 = "knee pads"
    52 = "compression sleeve youth baseball"
    53 = "tights for football"
    54 = "soccer pad"
    55 = "cycling pants for men"
    56 = "dry fit leggings men"
    57 = "compression calf sleeve men basketball"
    58 = "long compression shorts men"
    59 = "compression shorts long men"
    60 = "compression pants and tops for men"
    61 = "leggings knee length"
    62 = "mens softball gear"
    63 = "yoga after knee replacement"
    64 = "wrestling knee sleeve youth"
    65 = "tights compression"
    66 = "mens compression pants pack"
    67 = "boys running pants"
    68 = "knee pads thigh support"
    69 = "youth baseball pants black"
    70 = "knee pads biking adult"
    71 = "youth boys leggings"
    72 = "adult pants"
    73 = "youth baseball compression sleeves"
    74 = "calf sleeves for men football"
    75 = "padded knee sleeve"
    76 = "knee pad exercise"
    77 = "recovery pants men"
    78 = "mens tight"
    79 = "mens outdoor basketball"
    80 = "soccer pants youth"
    81 = "protective shorts"
    82 = "baseball hand guard"
    83 = "bee pants"
    84 = "mens protective pads"
    85 = "mens shorts above knee"
    86 = "basketball knee support for men"
    87 = "compression pads for surgery"
    88 = "snowboarding pants boys"
    89 = "basketball pants for girls"
    90 = "youth baseball pants girls"
    91 = "knee pads for exercise"
    92 = "girls compression knee sleeves"
    93 = "men yoga pant"
    94 = "yoga pants mens black"
    95 = "basketball aids"
    96 = "knee compression sleeve volleyball"
    97 = "impact advanced recovery"
    98 = "best basketball"
    99 = "baseball compression sleeve youth"
    100 = "compression knee sleeve padded"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 1).Value = $values[$row]
}
